$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AA2").Value = 32
$ws.Range("AB2").Value = 16
$ws.Range("AE2").Value = 22
$ws.Range("AF2").Value = 24
$ws.Range("AJ2").Value = 55
$ws.Range("AK2").Value = 32
$ws.Range("AN2").Value = 23
$ws.Range("AO2").Value = 14.5
$ws.Range("F2").Value = 3.1
$ws.Range("G2").Value = 3.2
$ws.Range("H2").Value = 2.34
$ws.Range("I2").Value = 2.36
$ws.Range("J2").Value = 3.8
$ws.Range("L2").Value = 1.35
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.24
$ws.Range("P2").Value = 2.3
$ws.Range("Q2").Value = 1.74
$ws.Range("R2").Value = 1.52
$ws.Range("S2").Value = 2.88
$ws.Range("T2").Value = 1.64
$ws.Range("V2").Value = 1.72
$ws.Range("W2").Value = 1.45
$ws.Range("Z2").Value = 16.5

# Row 3
$ws.Range("AB3").Value = 15
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 14
$ws.Range("AE3").Value = 32
$ws.Range("AG3").Value = 18.5
$ws.Range("AO3").Value = 24
$ws.Range("F3").Value = 3.3
$ws.Range("G3").Value = 3.85
$ws.Range("H3").Value = 2.16
$ws.Range("I3").Value = 2.36
$ws.Range("L3").Value = 1.44
$ws.Range("N3").Value = 3.7
$ws.Range("O3").Value = 1.35
$ws.Range("P3").Value = 1.89
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.34
$ws.Range("S3").Value = 3.55
$ws.Range("T3").Value = 1.83
$ws.Range("U3").Value = 2.02
$ws.Range("V3").Value = 1.73
$ws.Range("W3").Value = 1.36
$ws.Range("X3").Value = 16.5
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 17.5

# Row 4
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AF4").Value = 24
$ws.Range("AG4").Value = 1000
$ws.Range("F4").Value = 2.74
$ws.Range("H4").Value = 2.54
$ws.Range("I4").Value = 2.84
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.9
$ws.Range("L4").Value = 1.4
$ws.Range("N4").Value = 3.8
$ws.Range("Q4").Value = 1.9
$ws.Range("R4").Value = 1.37
$ws.Range("U4").Value = 2.16
$ws.Range("V4").Value = 1.56
$ws.Range("W4").Value = 1.48
$ws.Range("X4").Value = 18.5
$ws.Range("Y4").Value = 1000

# Row 5
$ws.Range("AC5").Value = 13
$ws.Range("AD5").Value = 65
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 55
$ws.Range("AL5").Value = 80
$ws.Range("I5").Value = 17.5
$ws.Range("J5").Value = 4.9
$ws.Range("K5").Value = 5.2
$ws.Range("L5").Value = 1.44
$ws.Range("N5").Value = 3.25
$ws.Range("P5").Value = 1.78
$ws.Range("R5").Value = 1.28
$ws.Range("S5").Value = 3.9
$ws.Range("T5").Value = 2.66
$ws.Range("U5").Value = 1.5
$ws.Range("V5").Value = 1.06

# Row 6
$ws.Range("AD6").Value = 19.5
$ws.Range("AE6").Value = 130
$ws.Range("AH6").Value = 20
$ws.Range("L6").Value = 1.38
$ws.Range("N6").Value = 4.2
$ws.Range("O6").Value = 1.27
$ws.Range("P6").Value = 2.08
$ws.Range("Q6").Value = 1.84
$ws.Range("R6").Value = 1.42
$ws.Range("S6").Value = 3.15
$ws.Range("T6").Value = 1.7
$ws.Range("U6").Value = 2.2
$ws.Range("Z6").Value = 90

# Row 7
$ws.Range("AH7").Value = 22
$ws.Range("AN7").Value = 5
$ws.Range("L7").Value = 1.25
$ws.Range("M7").Value = 1.02
$ws.Range("P7").Value = 2.84
$ws.Range("R7").Value = 1.74
$ws.Range("S7").Value = 2.18
$ws.Range("U7").Value = 2.3

# Row 8
$ws.Range("AB8").Value = 5.4
$ws.Range("AG8").Value = 11.5
$ws.Range("AH8").Value = 220
$ws.Range("AL8").Value = 460
$ws.Range("AN8").Value = 25
$ws.Range("F8").Value = 1.65
$ws.Range("G8").Value = 1.73
$ws.Range("I8").Value = 8.800000000000001
$ws.Range("J8").Value = 3.35
$ws.Range("L8").Value = 1.64
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 2.42
$ws.Range("O8").Value = 1.61
$ws.Range("P8").Value = 1.46
$ws.Range("Q8").Value = 2.9
$ws.Range("T8").Value = 2.6
$ws.Range("U8").Value = 1.51
$ws.Range("W8").Value = 2.36

# Row 9
$ws.Range("AD9").Value = 1000
$ws.Range("F9").Value = 1.8
$ws.Range("G9").Value = 1.87
$ws.Range("I9").Value = 5.8
$ws.Range("J9").Value = 3.65
$ws.Range("L9").Value = 1.47
$ws.Range("N9").Value = 3.35
$ws.Range("O9").Value = 1.37
$ws.Range("P9").Value = 1.78
$ws.Range("Q9").Value = 2.12
$ws.Range("R9").Value = 1.29
$ws.Range("S9").Value = 3.9
$ws.Range("T9").Value = 2
$ws.Range("U9").Value = 1.94
$ws.Range("W9").Value = 2.14

# Row 10
$ws.Range("AD10").Value = 990
$ws.Range("AK10").Value = 65
$ws.Range("G10").Value = 1.76
$ws.Range("H10").Value = 5.8
$ws.Range("I10").Value = 6.8
$ws.Range("J10").Value = 3.7
$ws.Range("L10").Value = 1.41
$ws.Range("N10").Value = 3.65
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 1.89
$ws.Range("Q10").Value = 1.98
$ws.Range("R10").Value = 1.34
$ws.Range("S10").Value = 3.7
$ws.Range("T10").Value = 1.9
$ws.Range("U10").Value = 1.87
$ws.Range("V10").Value = 1.19
$ws.Range("W10").Value = 2.3

# Row 11
$ws.Range("AB11").Value = 9.4
$ws.Range("AE11").Value = 200
$ws.Range("AF11").Value = 14.5
$ws.Range("AN11").Value = 21
$ws.Range("AO11").Value = 80
$ws.Range("G11").Value = 2.3
$ws.Range("H11").Value = 3.45
$ws.Range("I11").Value = 3.8
$ws.Range("L11").Value = 1.44
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 3.45
$ws.Range("O11").Value = 1.35
$ws.Range("P11").Value = 1.83
$ws.Range("Q11").Value = 2.06
$ws.Range("S11").Value = 3.65
$ws.Range("T11").Value = 1.78
$ws.Range("U11").Value = 2.02
$ws.Range("V11").Value = 1.36
$ws.Range("W11").Value = 1.76
$ws.Range("X11").Value = 16.5
$ws.Range("Y11").Value = 13.5

# Row 12
$ws.Range("AA12").Value = 70
$ws.Range("AC12").Value = 7.4
$ws.Range("AI12").Value = 290
$ws.Range("AJ12").Value = 36
$ws.Range("AL12").Value = 55
$ws.Range("AN12").Value = 34
$ws.Range("AO12").Value = 1000
$ws.Range("F12").Value = 2.4
$ws.Range("G12").Value = 2.48
$ws.Range("H12").Value = 3.45
$ws.Range("I12").Value = 3.55
$ws.Range("J12").Value = 3.25
$ws.Range("K12").Value = 3.4
$ws.Range("L12").Value = 1.54
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 3.05
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 1.67
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 1.26
$ws.Range("S12").Value = 4.7
$ws.Range("T12").Value = 1.96
$ws.Range("U12").Value = 1.92
$ws.Range("V12").Value = 1.39
$ws.Range("W12").Value = 1.68
$ws.Range("Y12").Value = 12.5

# Row 13
$ws.Range("AH13").Value = 60
$ws.Range("AJ13").Value = 90
$ws.Range("G13").Value = 2.08
$ws.Range("I13").Value = 4.7
$ws.Range("K13").Value = 3.6
$ws.Range("L13").Value = 1.46
$ws.Range("N13").Value = 3.4
$ws.Range("P13").Value = 1.79
$ws.Range("Q13").Value = 2.16
$ws.Range("R13").Value = 1.3
$ws.Range("S13").Value = 4
$ws.Range("W13").Value = 1.92

